# Auto-generated script applying the Anima_Profits market-data refresh diff.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# affected Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1624.125
$ws.Range("I70").Value = 1168
$ws.Range("J70").Value = 1897.8
$ws.Range("K70").Value = 3504
$ws.Range("L70").Value = 5693.4
$ws.Range("M70").Value = -3234
$ws.Range("N70").Value = -6233.4
$ws.Range("H73").Value = 1624.125
$ws.Range("I73").Value = 1168
$ws.Range("J73").Value = 1897.8
$ws.Range("K73").Value = 3504
$ws.Range("L73").Value = 5693.4
$ws.Range("M73").Value = -2568
$ws.Range("N73").Value = -7565.4
$ws.Range("H113").Value = 1902.4286
$ws.Range("J113").Value = 1953
$ws.Range("L113").Value = 1953
$ws.Range("N113").Value = -8461
$ws.Range("H132").Value = 4047.389
$ws.Range("I132").Value = 3573.1667
$ws.Range("J132").Value = 6418.5
$ws.Range("K132").Value = 10719.5001
$ws.Range("L132").Value = 19255.5
$ws.Range("M132").Value = -8189.500100000001
$ws.Range("N132").Value = -24315.5
$ws.Range("H137").Value = 1291.7
$ws.Range("I137").Value = 783
$ws.Range("K137").Value = 2349
$ws.Range("M137").Value = 201
$ws.Range("H141").Value = 3993.96
$ws.Range("I141").Value = 2102.647
$ws.Range("J141").Value = 8013
$ws.Range("K141").Value = 6307.941
$ws.Range("L141").Value = 24039
$ws.Range("M141").Value = -1127.941
$ws.Range("N141").Value = -34399

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 15000.111
$ws.Range("J23").Value = 11764.471
$ws.Range("L23").Value = 11764.471
$ws.Range("N23").Value = -12282.471
$ws.Range("H52").Value = 97379.336
$ws.Range("J52").Value = 97379.336
$ws.Range("L52").Value = 97379.336
$ws.Range("N52").Value = -98015.336
$ws.Range("H121").Value = 41245
$ws.Range("J121").Value = 41245
$ws.Range("L121").Value = 41245
$ws.Range("N121").Value = -44739
$ws.Range("H122").Value = 1854.5
$ws.Range("I122").Value = 1848
$ws.Range("J122").Value = 1861
$ws.Range("K122").Value = 5544
$ws.Range("L122").Value = 5583
$ws.Range("M122").Value = -3094
$ws.Range("N122").Value = -10483
$ws.Range("H132").Value = 5649.8486
$ws.Range("I132").Value = 5017.1924
$ws.Range("J132").Value = 7999.7144
$ws.Range("K132").Value = 15051.5772
$ws.Range("L132").Value = 23999.1432
$ws.Range("M132").Value = -12521.5772
$ws.Range("N132").Value = -29059.1432

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 14442.167
$ws.Range("I8").Value = 2663.25
$ws.Range("J8").Value = 38000
$ws.Range("K8").Value = 2663.25
$ws.Range("L8").Value = 38000
$ws.Range("M8").Value = -2523.25
$ws.Range("N8").Value = -38280
$ws.Range("H12").Value = 486.66666
$ws.Range("I12").Value = 486.66666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 486.66666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -318.66666
$ws.Range("N12").ClearContents()
$ws.Range("H134").Value = 3321.5557
$ws.Range("I134").Value = 3236.75
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 9710.25
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -7175.25
$ws.Range("N134").Value = -17070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 27333.334
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 36000
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 36000
$ws.Range("M13").Value = -9861
$ws.Range("N13").Value = -36278
$ws.Range("H31").Value = 1340.0962
$ws.Range("I31").Value = 1312.1538
$ws.Range("J31").Value = 1368.0385
$ws.Range("K31").Value = 1312.1538
$ws.Range("L31").Value = 1368.0385
$ws.Range("M31").Value = -1017.1538
$ws.Range("N31").Value = -1958.0385
$ws.Range("H34").Value = 1340.0962
$ws.Range("I34").Value = 1312.1538
$ws.Range("J34").Value = 1368.0385
$ws.Range("K34").Value = 1312.1538
$ws.Range("L34").Value = 1368.0385
$ws.Range("M34").Value = -1110.1538
$ws.Range("N34").Value = -1772.0385
$ws.Range("H122").Value = 1774.7273
$ws.Range("I122").Value = 1709
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5127
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -2677
$ws.Range("N122").Value = -10750
$ws.Range("H132").Value = 8335975.5
$ws.Range("I132").Value = 2394.7144
$ws.Range("K132").Value = 7184.1432
$ws.Range("M132").Value = -4654.1432

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 961.439
$ws.Range("I68").Value = 776.44446
$ws.Range("K68").Value = 2329.33338
$ws.Range("M68").Value = -1518.33338
$ws.Range("H71").Value = 961.439
$ws.Range("I71").Value = 776.44446
$ws.Range("K71").Value = 6988.00014
$ws.Range("M71").Value = -2932.00014
$ws.Range("H92").Value = 698.2778
$ws.Range("I92").Value = 653.75
$ws.Range("J92").Value = 733.9
$ws.Range("K92").Value = 1961.25
$ws.Range("L92").Value = 2201.7
$ws.Range("M92").Value = -713.25
$ws.Range("N92").Value = -4697.7
$ws.Range("H93").Value = 7357.143
$ws.Range("J93").Value = 8500
$ws.Range("L93").Value = 25500
$ws.Range("N93").Value = -29244
$ws.Range("H133").Value = 13031.3
$ws.Range("I133").Value = 766.6667
$ws.Range("J133").Value = 18287.572
$ws.Range("K133").Value = 2300.0001
$ws.Range("L133").Value = 54862.716
$ws.Range("M133").Value = 2759.9999
$ws.Range("N133").Value = -64982.716

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 35000
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -36118
$ws.Range("H132").Value = 3078.3
$ws.Range("I132").Value = 2684.1428
$ws.Range("J132").Value = 3998
$ws.Range("K132").Value = 8052.428400000001
$ws.Range("L132").Value = 11994
$ws.Range("M132").Value = -5522.428400000001
$ws.Range("N132").Value = -17054

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3625.6287
$ws.Range("I132").Value = 3090.8635
$ws.Range("J132").Value = 4530.615
$ws.Range("K132").Value = 9272.5905
$ws.Range("L132").Value = 13591.845
$ws.Range("M132").Value = -6742.5905
$ws.Range("N132").Value = -18651.845
$ws.Range("H136").Value = 1370.7646
$ws.Range("I136").Value = 1572.1818
$ws.Range("J136").Value = 1001.5
$ws.Range("K136").Value = 4716.5454
$ws.Range("L136").Value = 3004.5
$ws.Range("M136").Value = -2166.5454
$ws.Range("N136").Value = -8104.5

